# Add the new weekly ranking sheet for 2026-02-18
$wb = $excel.ActiveWorkbook

# Reference sheet that already has the highlighted 'same title released multiple
# volumes this week' style (fillId pointing at the light-yellow fgColor 00FFFACD),
# used below to copy that exact cell format onto the new sheet.
$refSheet = $wb.Worksheets.Item("2026-02-11")
$refStyledCell = $refSheet.Range("C2")

# Insert the new sheet after the last existing sheet and name it.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "2026-02-18"

# Header row
$ws.Range("A1").Value = "rank"
$ws.Range("B1").Value = "title"
$ws.Range("C1").Value = "volume"
$ws.Range("D1").Value = "publisher"

# Data rows: rank, title, volume (publisher column intentionally left blank,
# matching every other weekly sheet in this workbook).
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = '転生したらスライムだった件'
$ws.Cells.Item(2, 3).Value = 31
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 'チェンソーマン'
$ws.Cells.Item(3, 3).Value = 23
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 'BORUTO-ボルト- -TWO BLUE VORTEX-'
$ws.Cells.Item(4, 3).Value = 7
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = '魔入りました!入間くん'
$ws.Cells.Item(5, 3).Value = 47
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 'ミステリと言う勿れ'
$ws.Cells.Item(6, 3).Value = 16
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 'アオのハコ'
$ws.Cells.Item(7, 3).Value = 24
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = '転生したら第七王子だったので、気ままに魔術を極めます'
$ws.Cells.Item(8, 3).Value = 22
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = '信じていた仲間達にダンジョン奥地で殺されかけたがギフト『無限ガチャ』でレベル9999の仲間達を手に入れて元パーティーメンバーと世界に復讐&『ざまぁ!』します!'
$ws.Cells.Item(9, 3).Value = 21
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 'Sランクパーティから解雇された~『呪いのアイテム』しか作れませんが、その性能はアーティファクト級なり……!~'
$ws.Cells.Item(10, 3).Value = 13
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = '桃源暗鬼'
$ws.Cells.Item(11, 3).Value = 28
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = '聖者無双'
$ws.Cells.Item(12, 3).Value = 16
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = '刃牙らへん'
$ws.Cells.Item(13, 3).Value = 6
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = '極楽街'
$ws.Cells.Item(14, 3).Value = 6
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = 'レベル1から始まる召喚無双 THE COMIC'
$ws.Cells.Item(15, 3).Value = 3
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = 'パリピ孔明'
$ws.Cells.Item(16, 3).Value = 24
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = '転生したらスライムだった件 異聞 ~魔国暮らしのトリニティ~'
$ws.Cells.Item(17, 3).Value = 13
$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = 'レベル1から始まる召喚無双 THE COMIC'
$ws.Cells.Item(18, 3).Value = 2
$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = 'だれでも抱けるキミが好き'
$ws.Cells.Item(19, 3).Value = 8
$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = 'ダークギャザリング'
$ws.Cells.Item(20, 3).Value = 19
$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = '魔入りました!入間くん if Episode of 魔フィア'
$ws.Cells.Item(21, 3).Value = 7
$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = 'レベル1から始まる召喚無双 THE COMIC'
$ws.Cells.Item(22, 3).Value = 1
$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 2).Value = '転生したらスライムだった件 クレイマンREVENGE'
$ws.Cells.Item(23, 3).Value = 8
$ws.Cells.Item(24, 1).Value = 23
$ws.Cells.Item(24, 2).Value = '追放された転生王子、『自動製作』スキルで領地を爆速で開拓し最強の村を作ってしまう~最強クラフトスキルで始める、楽々領地開拓スローライフ~'
$ws.Cells.Item(24, 3).Value = 6
$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(25, 2).Value = '凶乱令嬢ニア・リストン 病弱令嬢に転生した神殺しの武人の華麗なる無双録'
$ws.Cells.Item(25, 3).Value = 8
$ws.Cells.Item(26, 1).Value = 25
$ws.Cells.Item(26, 2).Value = '100万の命の上に俺は立っている'
$ws.Cells.Item(26, 3).Value = 23
$ws.Cells.Item(27, 1).Value = 26
$ws.Cells.Item(27, 2).Value = '義妹にすべてを奪われたのに元婚約者(上司)が溺愛してきます。1'
$ws.Cells.Item(27, 3).Value = 1
$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = '辺境の薬師、都でSランク冒険者となる~英雄村の少年がチート薬で無自覚無双~'
$ws.Cells.Item(28, 3).Value = 11
$ws.Cells.Item(29, 1).Value = 28
$ws.Cells.Item(29, 2).Value = '2.5次元の誘惑'
$ws.Cells.Item(29, 3).Value = 25
$ws.Cells.Item(30, 1).Value = 29
$ws.Cells.Item(30, 2).Value = 'レベル1から始まる召喚無双 THE COMIC'
$ws.Cells.Item(30, 3).Value = 8
$ws.Cells.Item(31, 1).Value = 30
$ws.Cells.Item(31, 2).Value = 'ブルーロック'
$ws.Cells.Item(31, 3).Value = 37
$ws.Cells.Item(32, 1).Value = 31
$ws.Cells.Item(32, 2).Value = '僕の心のヤバイやつ'
$ws.Cells.Item(32, 3).Value = 13
$ws.Cells.Item(33, 1).Value = 32
$ws.Cells.Item(33, 2).Value = 'メダリスト'
$ws.Cells.Item(33, 3).Value = 14
$ws.Cells.Item(34, 1).Value = 33
$ws.Cells.Item(34, 2).Value = 'モンスターがあふれる世界になったので、好きに生きたいと思います'
$ws.Cells.Item(34, 3).Value = 14
$ws.Cells.Item(35, 1).Value = 34
$ws.Cells.Item(35, 2).Value = '社畜剣聖、配信者になる ~ブラックギルド会社員、うっかり会社用回線でS級モンスターを相手に無双するところを全国配信してしまう~(コミック)'
$ws.Cells.Item(35, 3).Value = 3
$ws.Cells.Item(36, 1).Value = 35
$ws.Cells.Item(36, 2).Value = '土かぶりのエレナ姫'
$ws.Cells.Item(36, 3).Value = 7
$ws.Cells.Item(37, 1).Value = 36
$ws.Cells.Item(37, 2).Value = '鵺の陰陽師'
$ws.Cells.Item(37, 3).Value = 13
$ws.Cells.Item(38, 1).Value = 37
$ws.Cells.Item(38, 2).Value = '傷だらけの公爵令嬢は、隣国の皇太子に溺愛される1'
$ws.Cells.Item(38, 3).Value = 1
$ws.Cells.Item(39, 1).Value = 38
$ws.Cells.Item(39, 2).Value = '無能は不要と言われ『時計使い』の僕は職人ギルドから追い出されるも、ダンジョンの深部で真の力に覚醒する THE COMIC'
$ws.Cells.Item(39, 3).Value = 2
$ws.Cells.Item(40, 1).Value = 39
$ws.Cells.Item(40, 2).Value = '傷モノの花嫁'
$ws.Cells.Item(40, 3).Value = 10
$ws.Cells.Item(41, 1).Value = 40
$ws.Cells.Item(41, 2).Value = 'ドローイング 最強漫画家はお絵描きスキルで異世界無双する!17'
$ws.Cells.Item(41, 3).Value = 17
$ws.Cells.Item(42, 1).Value = 41
$ws.Cells.Item(42, 2).Value = '税金で買った本'
$ws.Cells.Item(42, 3).Value = 18
$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 2).Value = '世界最強の魔女、始めました ~私だけ『攻略サイト』を見れる世界で自由に生きます~'
$ws.Cells.Item(43, 3).Value = 11
$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 2).Value = '葬送のフリーレン'
$ws.Cells.Item(44, 3).Value = 15
$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 2).Value = 'ハイスクールハックアンドスラッシュ'
$ws.Cells.Item(45, 3).Value = 1
$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 2).Value = '苔から始まる異世界ライフ'
$ws.Cells.Item(46, 3).Value = 1
$ws.Cells.Item(47, 1).Value = 46
$ws.Cells.Item(47, 2).Value = '育ちすぎたタマ'
$ws.Cells.Item(47, 3).Value = 1
$ws.Cells.Item(48, 1).Value = 47
$ws.Cells.Item(48, 2).Value = '「変なバイト見つけた」時給××万円の理由がヤバすぎる1'
$ws.Cells.Item(48, 3).Value = 1
$ws.Cells.Item(49, 1).Value = 48
$ws.Cells.Item(49, 2).Value = '世界で一番綺麗な姉はレベル1'
$ws.Cells.Item(49, 3).Value = 9
$ws.Cells.Item(50, 1).Value = 49
$ws.Cells.Item(50, 2).Value = '世界で一番綺麗な姉はレベル1'
$ws.Cells.Item(50, 3).Value = 1
$ws.Cells.Item(51, 1).Value = 50
$ws.Cells.Item(51, 2).Value = '神血の救世主~0.00000001%を引き当て最強へ~'
$ws.Cells.Item(51, 3).Value = 11
$ws.Cells.Item(52, 1).Value = 51
$ws.Cells.Item(52, 2).Value = '神血の救世主~0.00000001%を引き当て最強へ~'
$ws.Cells.Item(52, 3).Value = 12
$ws.Cells.Item(53, 1).Value = 52
$ws.Cells.Item(53, 2).Value = 'チート薬師のスローライフ'
$ws.Cells.Item(53, 3).Value = 14
$ws.Cells.Item(54, 1).Value = 53
$ws.Cells.Item(54, 2).Value = '転生したら平民でした。~生活水準に耐えられないので貴族を目指します~(コミック)'
$ws.Cells.Item(54, 3).Value = 7
$ws.Cells.Item(55, 1).Value = 54
$ws.Cells.Item(55, 2).Value = '金田一パパの事件簿'
$ws.Cells.Item(55, 3).Value = 3
$ws.Cells.Item(56, 1).Value = 55
$ws.Cells.Item(56, 2).Value = '29歳独身中堅冒険者の日常'
$ws.Cells.Item(56, 3).Value = 21
$ws.Cells.Item(57, 1).Value = 56
$ws.Cells.Item(57, 2).Value = '彼女、お借りします'
$ws.Cells.Item(57, 3).Value = 44
$ws.Cells.Item(58, 1).Value = 57
$ws.Cells.Item(58, 2).Value = '捨てられた地味王女は白狼殿下に溺愛される1'
$ws.Cells.Item(58, 3).Value = 1
$ws.Cells.Item(59, 1).Value = 58
$ws.Cells.Item(59, 2).Value = '二度目の人生では、お飾り王妃になりません!1'
$ws.Cells.Item(59, 3).Value = 1
$ws.Cells.Item(60, 1).Value = 59
$ws.Cells.Item(60, 2).Value = 'ハイスクールハックアンドスラッシュ'
$ws.Cells.Item(60, 3).Value = 2
$ws.Cells.Item(61, 1).Value = 60
$ws.Cells.Item(61, 2).Value = 'ハイスクールハックアンドスラッシュ'
$ws.Cells.Item(61, 3).Value = 3
$ws.Cells.Item(62, 1).Value = 61
$ws.Cells.Item(62, 2).Value = '苔から始まる異世界ライフ'
$ws.Cells.Item(62, 3).Value = 2
$ws.Cells.Item(63, 1).Value = 62
$ws.Cells.Item(63, 2).Value = '苔から始まる異世界ライフ'
$ws.Cells.Item(63, 3).Value = 3
$ws.Cells.Item(64, 1).Value = 63
$ws.Cells.Item(64, 2).Value = '三原ソフィアは怖ガール'
$ws.Cells.Item(64, 3).Value = 1
$ws.Cells.Item(65, 1).Value = 64
$ws.Cells.Item(65, 2).Value = '育ちすぎたタマ'
$ws.Cells.Item(65, 3).Value = 2
$ws.Cells.Item(66, 1).Value = 65
$ws.Cells.Item(66, 2).Value = '育ちすぎたタマ'
$ws.Cells.Item(66, 3).Value = 3
$ws.Cells.Item(67, 1).Value = 66
$ws.Cells.Item(67, 2).Value = 'ダンジョン教室'
$ws.Cells.Item(67, 3).Value = 1
$ws.Cells.Item(68, 1).Value = 67
$ws.Cells.Item(68, 2).Value = '16年目の復讐~奴らを地獄に送るまで1'
$ws.Cells.Item(68, 3).Value = 1
$ws.Cells.Item(69, 1).Value = 68
$ws.Cells.Item(69, 2).Value = '世界で一番綺麗な姉はレベル1'
$ws.Cells.Item(69, 3).Value = 8
$ws.Cells.Item(70, 1).Value = 69
$ws.Cells.Item(70, 2).Value = '世界で一番綺麗な姉はレベル1'
$ws.Cells.Item(70, 3).Value = 7
$ws.Cells.Item(71, 1).Value = 70
$ws.Cells.Item(71, 2).Value = '世界で一番綺麗な姉はレベル1'
$ws.Cells.Item(71, 3).Value = 6
$ws.Cells.Item(72, 1).Value = 71
$ws.Cells.Item(72, 2).Value = '世界で一番綺麗な姉はレベル1'
$ws.Cells.Item(72, 3).Value = 5
$ws.Cells.Item(73, 1).Value = 72
$ws.Cells.Item(73, 2).Value = '世界で一番綺麗な姉はレベル1'
$ws.Cells.Item(73, 3).Value = 4
$ws.Cells.Item(74, 1).Value = 73
$ws.Cells.Item(74, 2).Value = '世界で一番綺麗な姉はレベル1'
$ws.Cells.Item(74, 3).Value = 3
$ws.Cells.Item(75, 1).Value = 74
$ws.Cells.Item(75, 2).Value = '世界で一番綺麗な姉はレベル1'
$ws.Cells.Item(75, 3).Value = 2
$ws.Cells.Item(76, 1).Value = 75
$ws.Cells.Item(76, 2).Value = '絶対に死ぬモブ悪役令嬢は初恋がしたい 第1話'
$ws.Cells.Item(76, 3).Value = 1
$ws.Cells.Item(77, 1).Value = 76
$ws.Cells.Item(77, 2).Value = '超絶変身!! アースカイザー'
$ws.Cells.Item(77, 3).Value = 1
$ws.Cells.Item(78, 1).Value = 77
$ws.Cells.Item(78, 2).Value = '嫁入りのススメ~大正御曹司の強引な求婚~7'
$ws.Cells.Item(78, 3).Value = 7
$ws.Cells.Item(79, 1).Value = 78
$ws.Cells.Item(79, 2).Value = 'スーパーの裏でヤニ吸うふたり 通常版'
$ws.Cells.Item(79, 3).Value = 8
$ws.Cells.Item(80, 1).Value = 79
$ws.Cells.Item(80, 2).Value = 'レベル1から始まる召喚無双 THE COMIC'
$ws.Cells.Item(80, 3).Value = 6
$ws.Cells.Item(81, 1).Value = 80
$ws.Cells.Item(81, 2).Value = '呪術廻戦≡(モジュロ)'
$ws.Cells.Item(81, 3).Value = 1
$ws.Cells.Item(82, 1).Value = 81
$ws.Cells.Item(82, 2).Value = 'レベル1から始まる召喚無双 THE COMIC'
$ws.Cells.Item(82, 3).Value = 4
$ws.Cells.Item(83, 1).Value = 82
$ws.Cells.Item(83, 2).Value = 'レベル1から始まる召喚無双 THE COMIC'
$ws.Cells.Item(83, 3).Value = 5
$ws.Cells.Item(84, 1).Value = 83
$ws.Cells.Item(84, 2).Value = '峰子と魔王 ~異世界転移の若返り最強婆さん、最弱魔王と世直しする~'
$ws.Cells.Item(84, 3).Value = 1
$ws.Cells.Item(85, 1).Value = 84
$ws.Cells.Item(85, 2).Value = '中通りダイアリー'
$ws.Cells.Item(85, 3).Value = 1
$ws.Cells.Item(86, 1).Value = 85
$ws.Cells.Item(86, 2).Value = '後方見守り幼なじみが甘すぎる'
$ws.Cells.Item(86, 3).Value = 1
$ws.Cells.Item(87, 1).Value = 86
$ws.Cells.Item(87, 2).Value = 'おとずれナース ~精神科訪問看護とこころの記録~'
$ws.Cells.Item(87, 3).Value = 1
$ws.Cells.Item(88, 1).Value = 87
$ws.Cells.Item(88, 2).Value = '魔王城ホテルの悪役令嬢'
$ws.Cells.Item(88, 3).Value = 1
$ws.Cells.Item(89, 1).Value = 88
$ws.Cells.Item(89, 2).Value = '魔王城ホテルの悪役令嬢'
$ws.Cells.Item(89, 3).Value = 2
$ws.Cells.Item(90, 1).Value = 89
$ws.Cells.Item(90, 2).Value = '魔王城ホテルの悪役令嬢'
$ws.Cells.Item(90, 3).Value = 3
$ws.Cells.Item(91, 1).Value = 90
$ws.Cells.Item(91, 2).Value = '身代わり秒バレ令嬢の契約結婚なのに、騎士公爵が「絶対に離婚しない」と溺愛してくる'
$ws.Cells.Item(91, 3).Value = 1
$ws.Cells.Item(92, 1).Value = 91
$ws.Cells.Item(92, 2).Value = '身代わり秒バレ令嬢の契約結婚なのに、騎士公爵が「絶対に離婚しない」と溺愛してくる'
$ws.Cells.Item(92, 3).Value = 2
$ws.Cells.Item(93, 1).Value = 92
$ws.Cells.Item(93, 2).Value = '身代わり秒バレ令嬢の契約結婚なのに、騎士公爵が「絶対に離婚しない」と溺愛してくる'
$ws.Cells.Item(93, 3).Value = 3
$ws.Cells.Item(94, 1).Value = 93
$ws.Cells.Item(94, 2).Value = '悪役令嬢の遺言状1'
$ws.Cells.Item(94, 3).Value = 1
$ws.Cells.Item(95, 1).Value = 94
$ws.Cells.Item(95, 2).Value = '無能は不要と言われ『時計使い』の僕は職人ギルドから追い出されるも、ダンジョンの深部で真の力に覚醒する THE COMIC'
$ws.Cells.Item(95, 3).Value = 1
$ws.Cells.Item(96, 1).Value = 95
$ws.Cells.Item(96, 2).Value = '異世界迷宮でハーレムを'
$ws.Cells.Item(96, 3).Value = 1
$ws.Cells.Item(97, 1).Value = 96
$ws.Cells.Item(97, 2).Value = '三原ソフィアは怖ガール'
$ws.Cells.Item(97, 3).Value = 2
$ws.Cells.Item(98, 1).Value = 97
$ws.Cells.Item(98, 2).Value = '三原ソフィアは怖ガール'
$ws.Cells.Item(98, 3).Value = 3
$ws.Cells.Item(99, 1).Value = 98
$ws.Cells.Item(99, 2).Value = 'ダンジョン教室'
$ws.Cells.Item(99, 3).Value = 2
$ws.Cells.Item(100, 1).Value = 99
$ws.Cells.Item(100, 2).Value = 'ダンジョン教室'
$ws.Cells.Item(100, 3).Value = 3
$ws.Cells.Item(101, 1).Value = 100
$ws.Cells.Item(101, 2).Value = 'サベージ・キャッスル~堕落の迷宮~ 第1話'
$ws.Cells.Item(101, 3).Value = 1

# Re-apply the highlight fill to the 'volume' cells for titles that had more than
# one volume published the same week, reusing the workbook's existing style so we
# don't create a duplicate/near-duplicate style entry.
$refStyledCell.Copy()
$ws.Range("C15").PasteSpecial(-4122)
$refStyledCell.Copy()
$ws.Range("C18").PasteSpecial(-4122)
$refStyledCell.Copy()
$ws.Range("C22").PasteSpecial(-4122)
$refStyledCell.Copy()
$ws.Range("C27").PasteSpecial(-4122)
$refStyledCell.Copy()
$ws.Range("C35").PasteSpecial(-4122)
$refStyledCell.Copy()
$ws.Range("C38:C39").PasteSpecial(-4122)
$refStyledCell.Copy()
$ws.Range("C45:C48").PasteSpecial(-4122)
$refStyledCell.Copy()
$ws.Range("C50").PasteSpecial(-4122)
$refStyledCell.Copy()
$ws.Range("C55").PasteSpecial(-4122)
$refStyledCell.Copy()
$ws.Range("C58:C68").PasteSpecial(-4122)
$refStyledCell.Copy()
$ws.Range("C74:C77").PasteSpecial(-4122)
$refStyledCell.Copy()
$ws.Range("C81").PasteSpecial(-4122)
$refStyledCell.Copy()
$ws.Range("C84:C101").PasteSpecial(-4122)

$excel.CutCopyMode = 0
